# Apply crypto price/volume updates from the latest GitHub Actions scrape.
# Source data: coinranking.com snapshot refreshed on Wed May 31 17:32:41 UTC 2023.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Coin / Link / Price / Volume(1h) figures for every affected row.
# A leading apostrophe forces the Price (column D) cells to stay text, matching
# the workbook convention where values such as "1.000" / "305.80" are stored as
# strings (not numbers) -- same as the pre-existing cells in this column.

$ws.Range("D2").Value = "'26.958.16"
$ws.Range("E2").Value = "  -2.51%  "

$ws.Range("D3").Value = "'1.861.03"
$ws.Range("E3").Value = "  -2.04%  "

$ws.Range("D4").Value = "'0.9992"
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'305.88"
$ws.Range("E5").Value = "  -1.98%  "

$ws.Range("D6").Value = "'0.9993"
$ws.Range("E6").Value = "  -0.05%  "

$ws.Range("D7").Value = "'0.5051"
$ws.Range("E7").Value = "  -3.51%  "

$ws.Range("D8").Value = "'0.3740"
$ws.Range("E8").Value = "  -0.85%  "

$ws.Range("D9").Value = "'0.07138"
$ws.Range("E9").Value = "  -1.27%  "

$ws.Range("D10").Value = "'0.8821"
$ws.Range("E10").Value = "  -1.17%  "

$ws.Range("D11").Value = "'20.60"
$ws.Range("E11").Value = "  -2.12%  "

$ws.Range("B12").Value = "WrappedEther"  # row 12 coin name swap
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"  # row 12 link swap
$ws.Range("D12").Value = "'1.885.81"
$ws.Range("E12").Value = "  -0.52%  "

$ws.Range("B13").Value = "TRON"  # row 13 coin name swap
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"  # row 13 link swap
$ws.Range("D13").Value = "'0.07565"
$ws.Range("E13").Value = "  -0.69%  "

$ws.Range("D14").Value = "'5.292"
$ws.Range("E14").Value = "  -2.63%  "

$ws.Range("D15").Value = "'89.08"
$ws.Range("E15").Value = "  -3.01%  "

$ws.Range("D16").Value = "'0.9993"
$ws.Range("E16").Value = "  -0.07%  "

$ws.Range("D17").Value = "'0.000008410"
$ws.Range("E17").Value = "  -3.32%  "

$ws.Range("D18").Value = "'14.09"
$ws.Range("E18").Value = "  -2.30%  "

$ws.Range("D19").Value = "'0.9997"
$ws.Range("E19").Value = "  +0.03%  "

$ws.Range("D20").Value = "'26.996.35"
$ws.Range("E20").Value = "  -2.51%  "

$ws.Range("D21").Value = "'5.030"
$ws.Range("E21").Value = "  -1.79%  "

$ws.Range("D22").Value = "'2.100.12"
$ws.Range("E22").Value = "  -1.71%  "

$ws.Range("D23").Value = "'10.47"
$ws.Range("E23").Value = "  -3.11%  "

$ws.Range("D24").Value = "'6.454"
$ws.Range("E24").Value = "  -1.67%  "

$ws.Range("E25").Value = "  -1.16%  "

$ws.Range("D26").Value = "'147.18"
$ws.Range("E26").Value = "  -3.83%  "

$ws.Range("D27").Value = "'17.94"
$ws.Range("E27").Value = "  -1.73%  "

$ws.Range("D28").Value = "'2.093"
$ws.Range("E28").Value = "  -2.90%  "

$ws.Range("D29").Value = "'112.59"
$ws.Range("E29").Value = "  -1.52%  "

$ws.Range("D30").Value = "'4.663"
$ws.Range("E30").Value = "  -3.34%  "

$ws.Range("D31").Value = "'4.693"
$ws.Range("E31").Value = "  -2.69%  "

$ws.Range("D32").Value = "'0.09038"

$ws.Range("D33").Value = "'0.05130"
$ws.Range("E33").Value = "  -2.50%  "

$ws.Range("D34").Value = "'3.023"
$ws.Range("E34").Value = "  -4.58%  "

$ws.Range("D35").Value = "'1.150"
$ws.Range("E35").Value = "  -6.70%  "

$ws.Range("D36").Value = "'0.7243"
$ws.Range("E36").Value = "  -5.52%  "

$ws.Range("D37").Value = "'0.02037"
$ws.Range("E37").Value = "  -1.69%  "

$ws.Range("D38").Value = "'3.036"
$ws.Range("E38").Value = "  -0.59%  "

$ws.Range("D39").Value = "'2.459"
$ws.Range("E39").Value = "  -5.24%  "

$ws.Range("D40").Value = "'1.077"
$ws.Range("E40").Value = "  -1.18%  "

$ws.Range("E41").Value = "  -3.69%  "

$ws.Range("D42").Value = "'6.516"
$ws.Range("E42").Value = "  -1.56%  "

$ws.Range("D43").Value = "'115.28"
$ws.Range("E43").Value = "  +1.97%  "

$ws.Range("D44").Value = "'8.242"
$ws.Range("E44").Value = "  -2.14%  "

$ws.Range("E45").Value = "  -2.55%  "

$ws.Range("D46").Value = "'0.9988"
$ws.Range("E46").Value = "  -0.06%  "

$ws.Range("D47").Value = "'0.4601"
$ws.Range("E47").Value = "  -3.40%  "

$ws.Range("D48").Value = "'9.973"
$ws.Range("E48").Value = "  -3.85%  "

$ws.Range("D49").Value = "'1.564"
$ws.Range("E49").Value = "  -2.77%  "

$ws.Range("D50").Value = "'36.50"
$ws.Range("E50").Value = "  -0.50%  "

$ws.Range("D51").Value = "'63.87"
$ws.Range("E51").Value = "  -3.64%  "
